# Review_322.docx edit: update title/date, rewrite the body paragraphs for the
# new paper (EFFICIENT REINFORCEMENT LEARNING WITH LARGE LANGUAGE MODEL PRIORS),
# and drop the two trailing paragraphs (old sign-off + old arxiv link), since the
# new arxiv link now lives in the paragraph that used to hold the last body text.

$d = $word.ActiveDocument

$d.Content.Find.Execute("⚡️🚀המאמר היומי של מייק -17.10.24: ⚡️🚀", $true, $false, $false, $false, $false, $true, 1, $false, "⚡️🚀המאמר היומי של מייק -16.10.24: ⚡️🚀", 2) | Out-Null
$d.Content.Find.Execute("EQUIVARIANT CONTRASTIVE LEARNING", $true, $false, $false, $false, $false, $true, 1, $false, "EFFICIENT REINFORCEMENT LEARNING WITH LARGE LANGUAGE MODEL PRIORS", 2) | Out-Null
$d.Content.Find.Execute("היום נסקור מאמר שפורסם לפני שנתיים וחצי בנושא למידה ניגודית (contrastive learning). הנושא עצמו תמיד עניין אותי וסקרתי לא מעט מאמרים אבל חייב להגיד שבזמן האחרון שטף המאמרים על CL די נחלש. כאמור המאמר הזה שראה אור לפני שנתיים מציע שכלול לשיטה הקלאסית לבנייה של ייצוג דאטה (אמבדינג) באמצעות CL.", $true, $false, $false, $false, $false, $true, 1, $false, "היום נסקור מאמר שהוא נראה די כבד מתמטית (הרבה נוסחאות ומלל שנראה מתמטי) אבל הרעיון מאחוריו הוא די פשוט וקל להסבר. אנחנו אוהבים למנף  את עוצמתם של מודלי שפה למשימות רבות (ולא תמיד לכאלו שהם מסוגלים לבצע כמו שצריך לפחות כרגע).", 2) | Out-Null
$d.Content.Find.Execute("בגדול CL היא שיטה לבניית ייצוג של דאטה כאשר העיקרון המוביל הוא לקרב ייצוגי פיסות דאטה דומות(זוגות חיוביים) ולהרחיק ייצוגים של פיסות דאטה לא דומות (שליליים).  זוגות דוגמאות חיוביים (במקרה של דאטה לא מתויג) נבחרות כאוגמנטציות שונות של דוגמא (עבור תמונות זה יכול להיות הזזה, סיבוב וכדומה) ואילו זוגות השליליים נבחרים באקראי מהדאטהסט. ", $true, $false, $false, $false, $false, $true, 1, $false, "המאמר מציע להשתמש במודל שפה כפריור עבור סוכנים במשימות בהם הם צריכים לבצע SDM או sequential decision making. המאמר נותן בתור דוגמא משחק overcooked כאשר הסוכן צריך לבצע משימות בישול שונות בהתבסס על מצב המטבח שבו הוא מבשל אותם. המטרה של הסוכן היא לחזות את הפעולה הבא (באמצעות תיאור טקסטואלי) כאשר התגמול הוא ביצוע נכון של המשימה (הכנה של מנה לפי המתכון :)).", 2) | Out-Null
$d.Content.Find.Execute("אולם יש לא מעט בעיות עם הגישה הזו הקשורות לבחירת זוגות של דוגמאות חיוביות - למשל שני פאצ'ים באותה התמונה עלולים להכיל תוכן סמנטי שונה שלא נרצה לקרב את ייצוגיהם (הוצעו מספר פתרונות לסוגיה זו בעבר וחלקן סקרתי). בנוסף אולי היינו רוצים לקבל ייצוגים שונים (ולא מאוד קרובים) של טרנספורמציות מסוימות של אותה התמונה (נגיד סיבוב או הזזה) למשימת downstream ספציפית. ", $true, $false, $false, $false, $false, $true, 1, $false, "כאמור המטרה כאן היא לחזות את הפעולה הבאה עבור הסוכן (המתוארת) על ידי הטקסט כאשר המצב (state) גם מתואר על ידי טקסט. בגדול מאוד אנו מתחילים ממודל אחד (הפריור P) עבור חיזוי המצב הבא (מהמצב הקודם והפעולה) ועבור חיזוי הפעולה הבאה בהינתן המצב (מתואר על ידי התפלגות Q_h). המטרה כאן היא ללמוד את Q_h כאשר ממקסמת התגמול הצפוי ושומרת את התפלגות Q קרובה לפריור P (זוכרים PPO שהתפרסם מאוד לפני שנתיים כאשר OpenAI השתמשו בו ל-RLHF לאימון מודלי שפה). המרחק כמובן ניתן על ידי ה-KL 🙂", 2) | Out-Null
$d.Content.Find.Execute("כלומר היינו רוצים להשרות יחס נתון T_i בין ייצוגי התמונה ההתחלתית I ולייצוג התמונה אחרי טרנספורמציה T (נקרא לה I_T). כלומר אנו רוצים לבנות ייצוג p כך ש:", $true, $false, $false, $false, $false, $true, 1, $false, "אז הפעולה הבאה a_t (כלומר גנרוט התיאור הטקסטואלי שלה) מתבצע באופן הבא. דוגמים כמה גרסאות של a_t עם P מחשבים את הנראות שלהם לפי Q הנלמד, מנרמלים עם הסופטמקס ודוגמים את הפעולה הבאה כאשר מטרת התהליך מקסום של התגמול הצפוי (עם הרגולריזציה שהסברנו עליה קודם).", 2) | Out-Null
$d.Content.Find.Execute("p(T(I)) = I_T(p(I))", $true, $false, $false, $false, $false, $true, 1, $false, "כמובן שניתן לעשות את זה בכמה אופנים: בצורה של online דרך שערוך של פונקציית Q של הזוג (מצב, פעולה) כאשר פונקציית Q קשורה להתפלגות Q_h של הפעולה הבא שנידונה בפסקה הקודמת (עניין של נרמול נכון). ניתן לעשות את זה גם באמצעות offline עם איזה פוליסי טוב ידוע של המומחים כאשר המטרה היא גם שערוך של פונקציית Q שבאמצעותה ניתן לשערך (לקבל) את Q_h עבור חיזוי הפעולה הבא. ניתן לעשות את זה גם באמצעות שיטה דומה ל-PPO אבל בכל המקרים הפריור הוא ההתפלגות המושרית על ידי מודל שפה נתון.", 2) | Out-Null
$d.Content.Find.Execute("וזה בדיוק מה שנקרא equivariance. למעשה CL הסטנדרטי הוא מקרה פרטי של equivariance שעבורן T_i הינה טרנספורמצית זהות וזה נקרא אינווריאנטיות של הייצוג תחת טרנספורמציית T.", $true, $false, $false, $false, $false, $true, 1, $false, "מאמר מעניין בקיצור…", 2) | Out-Null
$d.Content.Find.Execute(" וזה בדיוק מה שהמאמר עושה. למעשה המחברים מציעים לאמן ייצוג ששומר על אינווריאנטיות עבור טרנספורמציות מסוימות (כמו בCL הסטנדרטי) ו אוכף בנוסף equivariance מוגדר לטרנספורמציות מקבוצה נתונה G המתאימה למשימת downstream שיש לנו ביד. כלומר לכל טרנספורמציה מ-G אנו מגדירים מראש את הטרנספורמציה ה-equivariant שלה (שיכולה להיות חברה ב-G גם כן) ומאמנים את הייצוג כך שהיחס ה-equivariance ביניהם יתקיים. מבחינה פרקטית הלוס הוא סכום משוקלל של הלוסים של CL הסטנדרטי ו ה-ECL. ", $true, $false, $false, $false, $false, $true, 1, $false, "https://arxiv.org/pdf/2410.07927", 2) | Out-Null

# Drop the two now-obsolete trailing paragraphs (old sign-off line and old arxiv
# link). Walk from the last paragraph backwards so deleting one does not shift the
# index of a paragraph we still need to check/delete.
$obsolete = @("מאמר חמוד - מחר או היום בערב אסקור את מאמר ההמשך שלו…", "https://arxiv.org/abs/2111.00899")
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $t = $d.Paragraphs($i).Range.Text
    foreach ($o in $obsolete) {
        if ($t -like "*$o*") {
            $d.Paragraphs($i).Range.Delete()
            break
        }
    }
}

Write-Host "Final paragraph count:" $d.Paragraphs.Count
